# Update the dSF (column F) values per the "repull data, push all data,
# mean calculation" pass. Rows 6, 21 and 25 were not touched by this pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = -1
    4  = 3
    5  = 2
    7  = -1
    8  = -5
    9  = -4
    10 = -6
    11 = -6
    12 = 2
    13 = -3
    14 = 2
    15 = -2
    16 = -3
    17 = 3
    18 = -2
    19 = -3
    20 = -1
    22 = 10
    23 = -3
    24 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
